# Upgrade "Fund Transaction Statement" sheet to support cash dividend per-share
# reporting: insert a new column N ("现金分红" / "{.dividendAmountPerShare}")
# between "成本均价" (col M) and "基金公司" (col N -> now O), shifting the
# remaining header/data columns one position to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fund Transaction Statement")

# Insert a new column before column N; Excel automatically shifts the
# existing N:S columns (and their formatting/merged ranges) to O:T.
$ws.Columns("N:N").Insert()

# Populate the two new header cells that the inserted column introduces.
$ws.Range("N2").Value = "现金分红"
$ws.Range("N3").Value = "{.dividendAmountPerShare}"
